$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 412; this shifts the existing rows 412:518
# down to 413:519 and extends the used range to A1:R519.
$ws.Rows.Item(412).Insert()

# Populate the newly inserted row 412 with the new "Repollo" record.
$ws.Range("A412").Value = 5
$ws.Range("B412").Value = "Macroferia Regional de Talca"
$ws.Range("C412").Value = "Maule"
$ws.Range("D412").Value = 45135
$ws.Range("D412").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E412").Value = 7
$ws.Range("F412").Value = 100112006
$ws.Range("G412").Value = "Repollo"
$ws.Range("H412").Value = "Crespo record"
$ws.Range("I412").Value = "Primera"
$ws.Range("J412").Value = 5000
$ws.Range("K412").Value = 600
$ws.Range("L412").Value = 600
$ws.Range("M412").Value = 600
$ws.Range("N412").Value = "`$/unidad"
$ws.Range("O412").Value = "Región del Maule"
$ws.Range("P412").Value = 600
$ws.Range("Q412").Value = 1
$ws.Range("R412").Value = "Hortaliza"
